# Update MSME definitions: replace the literal "<br/>" markers inside the
# "Number of employees" and "Turnover" definition cells (for Medium and
# Large enterprise rows) with real line breaks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

$cells = @("B21", "D21", "B22", "D22")
foreach ($addr in $cells) {
    $rng = $ws.Range($addr)
    $text = $rng.Value2
    $rng.Value2 = $text.Replace("<br/>", $nl)
}
